$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A82").Value = "select * from message_table;"
$ws.Range("A83").Value = "SELECT ""id"", ""customer_number"", ""customer_name"", ""client_id"", ""client_name"", ""captain_name"", ""delivery_date"", ""order_status"", ""shop_to_delivery_km"", ""order_created_at"", ""order_accepted_at"", ""start_ride_at"", ""reached_shop_at"", ""order_picked_at"", ""shipped_at"", ""reached_dest_at"", ""final_status_at"", ""cancellation_reason"" FROM ""updated_table"" WHERE ""client_name"" = 'Dominos Pizza' AND ""order_created_at"" >= '2025-09-19 15:20:02.178437' - INTERVAL '1 day' AND ""order_created_at"" <= '2025-09-19 15:20:02.178437' LIMIT 10;"
$ws.Range("A84").Value = "SELECT ""id"", ""customer_number"", ""customer_name"", ""order_status"", ""delivery_date"", ""final_status_at"" FROM updated_table WHERE ""client_name"" = 'Dominos Pizza' LIMIT 10;"
